# Weekly fruit/vegetable price update: a new weekly record for "Ajo"
# (garlic) was added to the Mercado Mayorista Lo Valledor de Santiago
# series. In the source system the new observation is inserted right
# before the existing row for the same product family (row 327),
# pushing all subsequent rows (327-345) down by one (to 328-346).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 327; rows 327:345 shift down to 328:346.
$ws.Rows("327:327").Insert()

# Populate the newly inserted row with the latest weekly record.
$ws.Range("A327").Value = 6
$ws.Range("B327").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C327").Value = "Metropolitana"
$ws.Range("D327").Value = 44461
$ws.Range("E327").Value = 13
$ws.Range("F327").Value = 100112003
$ws.Range("G327").Value = "Ajo"
$ws.Range("H327").Value = "Chino"
$ws.Range("I327").Value = "Primera"
$ws.Range("J327").Value = 2600
$ws.Range("K327").Value = 14500
$ws.Range("L327").Value = 15000
$ws.Range("M327").Value = 14788
$ws.Range("N327").Value = "`$/caja 10 kilos"
$ws.Range("O327").Value = "China"
$ws.Range("P327").Value = 1479
$ws.Range("Q327").Value = 10
$ws.Range("R327").Value = "Hortaliza"
